# LOM3229.xlsx -- fix misaligned shared-string references introduced when
# "Docentes responsaveis" (2 rows) were inserted without shifting the rows
# below; also adds the new Programa resumido / Programa / Bibliografia text
# and re-splits the merged column A/B width definition.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Clear cells that must become blank in the final layout ---
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
$ws.Range("B19").Clear()
$ws.Range("C19").Clear()
$ws.Range("B24").Clear()
$ws.Range("C24").Clear()

# --- 2) Update cell values that already existed (kept current style) ---
$ws.Range("B10").Value = 'Apresentar as técnicas experimentais de preparação materialográfica e de caracterização de materiais.'
$ws.Range("C10").Value = 'Apresentar as técnicas experimentais de preparação materialográfica e de caracterização de materiais.'
$ws.Range("B13").Value = '6495737 - Durval Rodrigues Junior'
$ws.Range("C13").Value = '6495737 - Durval Rodrigues Junior'
$ws.Range("B14").Value = '1643715 - Paulo Atsushi Suzuki'
$ws.Range("C14").Value = '1643715 - Paulo Atsushi Suzuki'
$ws.Range("A15").Value = 'Programa resumido:'
$ws.Range("B15").Value = 'Difração de raios X. Materialografia. Microscopia óptica. Microscopia eletrônica. Análise térmica.'
$ws.Range("C15").Value = 'Difração de raios X. Materialografia. Microscopia óptica. Microscopia eletrônica. Análise térmica.'
$ws.Range("A16").Value = 'Short syllabus:'
$ws.Range("B16").Value = 'X-ray diffraction. Materialography. Optical microscopy. Electron microscopy. Thermal analysis.'
$ws.Range("C16").Value = 'X-ray diffraction. Materialography. Optical microscopy. Electron microscopy. Thermal analysis.'
$ws.Range("A17").Value = 'Programa:'
$ws.Range("A18").Value = 'Syllabus:'
$ws.Range("B18").Value = 'The microstructure of materials. Crystal lattices and systems, space groups and symmetry, most common types of crystal structures. Stereographic projection. Direction of the diffracted beam and Bragg''s law. Intensity of the diffracted beam. Methods of X-ray diffraction.Materialographic sample preparation: cutting, embedding, sanding and polishing. Chemical etching techniques to reveal phases. Fundamentals of quantitative materialography. Optical microscopy. Electron microscopy techniques: scanning and transmission. Chemical analysis of microregions: energy dispersive spectroscopy. Thermal analysis techniques: differential thermal analysis, differential scanning calorimetry and thermogravimetric analysis.'
$ws.Range("C18").Value = 'The microstructure of materials. Crystal lattices and systems, space groups and symmetry, most common types of crystal structures. Stereographic projection. Direction of the diffracted beam and Bragg''s law. Intensity of the diffracted beam. Methods of X-ray diffraction.Materialographic sample preparation: cutting, embedding, sanding and polishing. Chemical etching techniques to reveal phases. Fundamentals of quantitative materialography. Optical microscopy. Electron microscopy techniques: scanning and transmission. Chemical analysis of microregions: energy dispersive spectroscopy. Thermal analysis techniques: differential thermal analysis, differential scanning calorimetry and thermogravimetric analysis.'
$ws.Range("A19").Value = 'Avaliação:'
$ws.Range("A20").Value = 'Método:'
$ws.Range("B20").Value = 'Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo.'
$ws.Range("C20").Value = 'Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo.'
$ws.Range("A21").Value = 'Critério:'
$ws.Range("B21").Value = 'Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3'
$ws.Range("C21").Value = 'Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3'
$ws.Range("A22").Value = 'Norma de recuperação:'
$ws.Range("B23").Value = 'PADILHA, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985.
MURPHY, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001.
WU, Q.; MERCHANT, F.; CASTLEMAN, K. Microscope Image Processing, Academic Press, 2008.
CULLITY, B. D.; STOCK, S. R. Elements of X-Ray Diffraction, Prentice Hall, 2001.
YACOBI, B. G.; HOLT, D. B.; KAZMERSKI, L. L. Microanalysis of Solids. Plenum Press, New York, 1994.
HATAKEYAMA, T.; ZHENHAI, L. Handbook of Thermal Analysis, Wiley, 1999.
HAINES, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002.'
$ws.Range("C23").Value = 'PADILHA, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985.
MURPHY, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001.
WU, Q.; MERCHANT, F.; CASTLEMAN, K. Microscope Image Processing, Academic Press, 2008.
CULLITY, B. D.; STOCK, S. R. Elements of X-Ray Diffraction, Prentice Hall, 2001.
YACOBI, B. G.; HOLT, D. B.; KAZMERSKI, L. L. Microanalysis of Solids. Plenum Press, New York, 1994.
HATAKEYAMA, T.; ZHENHAI, L. Handbook of Thermal Analysis, Wiley, 1999.
HAINES, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002.'
$ws.Range("B25").Value = 'LOB1021 -  Física IV  (Requisito)
'
$ws.Range("C25").Value = 'LOB1021 -  Física IV  (Requisito)
'

# --- 3) Write brand-new cells, then clone the format from column template cells ---
$ws.Range("B17").Value = 'A microestrutura dos materiais. Sistemas e reticulados cristalinos, grupos espaciais e simetria, tipos mais comuns de estruturas cristalinas. Projeção estereográfica. Direção do feixe difratado e a lei de Bragg. Intensidade do feixe difratado. Métodos de difração de raios X. Preparação materialográfica de amostras: corte, embutimento, lixamento e polimento. Técnicas de ataque químico para revelação de fases. Fundamentos de materialografia quantitativa. Microscopia óptica. Técnicas de microscopia eletrônica: varredura e transmissão. Análise química de microrregiões: espectroscopia de energia dispersiva. Técnicas de análise térmica: análise térmica diferencial, calorimetria exploratória diferencial e análise termogravimétrica.'
$ws.Range("C17").Value = 'A microestrutura dos materiais. Sistemas e reticulados cristalinos, grupos espaciais e simetria, tipos mais comuns de estruturas cristalinas. Projeção estereográfica. Direção do feixe difratado e a lei de Bragg. Intensidade do feixe difratado. Métodos de difração de raios X. Preparação materialográfica de amostras: corte, embutimento, lixamento e polimento. Técnicas de ataque químico para revelação de fases. Fundamentos de materialografia quantitativa. Microscopia óptica. Técnicas de microscopia eletrônica: varredura e transmissão. Análise química de microrregiões: espectroscopia de energia dispersiva. Técnicas de análise térmica: análise térmica diferencial, calorimetria exploratória diferencial e análise termogravimétrica.'
$ws.Range("B22").Value = 'Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C22").Value = 'Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("A23").Value = 'Bibliografia:'
$ws.Range("A24").Value = 'Requisitos:'
$ws.Range("B26").Value = 'LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)
'
$ws.Range("C26").Value = 'LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)
'
$ws.Range("B27").Value = 'LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)
'
$ws.Range("C27").Value = 'LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)
'

$ws.Range("A3").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("B27").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4) Row heights ---
$ws.Rows.Item(13).EntireRow.AutoFit()
$ws.Rows.Item(14).EntireRow.AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).EntireRow.AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120
$ws.Rows.Item(24).EntireRow.AutoFit()
$ws.Rows.Item(26).RowHeight = 30
$ws.Rows.Item(27).RowHeight = 30

# --- 5) Dimension grows to C27; re-split the col A/B width definition ---
# (column B already carries its own 60.7109375-wide <col> override further
#  along, so nudging column B re-splits the merged "min=1 max=2" range into
#  its own "min=1 max=1" entry without touching column As stored width)
$ws.Columns.Item(2).ColumnWidth = 60.7109375
